$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) / Volume(1h) (E) columns with the latest scrape values.
# These columns are stored as text in the workbook (thousands-dot-separated
# prices like "51.060.63", percentage strings padded with spaces, and plain
# decimals where a trailing zero is significant, e.g. "7.60"). Force Text
# format on the Price cells before assigning so Excel keeps them as literal
# strings instead of silently reinterpreting them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.060.63"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.943.41"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "376.54"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.27"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.409.38"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.14"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.941.74"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.996"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.01"
$ws.Range("E18").Value = "  +47.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.054.88"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -7.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.36"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.78"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +7.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.13"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.46"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.62"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  -6.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.00"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.74"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.06"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.46"
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.46"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.33"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.29"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.272"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.987.80"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0328"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("E51").Value = "  +2.20%  "
